# The scraper dropped the "reviews_count" field, so that column (E) is
# removed entirely and every field to its right shifts one column to
# the left: reviews_average F->E, latitude G->F, longitude H->G,
# is_permanently_closed I->H, gmaps_link J->I, latest_review_date K->J.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("E:E").Delete()
